$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: 57576 -> 58745
$ws.Range("B2").Value = 58745

# C2: remove the "Ovaliderad" cell entirely
$ws.Range("C2").ClearContents()

# J2 / N2 / AF2: new empty (inline string) cells
$ws.Range("J2").Value = "'"
$ws.Range("J2").Style = "Normal"
$ws.Range("N2").Value = "'"
$ws.Range("N2").Style = "Normal"
$ws.Range("AF2").Value = "'"
$ws.Range("AF2").Style = "Normal"

# Q2 / R2: truncate to whole numbers
$ws.Range("Q2").Value = 559941
$ws.Range("R2").Value = 6257482

# Z2 and AB2: remove the "00:00" cells entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# AE2: TRUE -> FALSE
$ws.Range("AE2").Value = $false
